$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich-text runs) ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Cells that change between numeric and placeholder-text representation ---
# (copy format+value from a same-shaped donor cell elsewhere on the sheet,
#  then overwrite with the real target value where the donor differs)
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("L14").Copy($ws.Range("E14"))
$ws.Range("H14").Copy($ws.Range("M15"))
$ws.Range("M15").Value = 0
$ws.Range("G14").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("G14").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("H14").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 0
$ws.Range("G14").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1
$ws.Range("G14").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 4
$ws.Range("H14").Copy($ws.Range("E17"))
$ws.Range("E17").Value = -75
$ws.Range("G14").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 3
$ws.Range("H14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -66.666666666666
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("L14").Copy($ws.Range("E20"))
$ws.Range("C14").Copy($ws.Range("F20"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("L14").Copy($ws.Range("E28"))

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -75
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -57.142857142857
$ws.Range("L16").Value = -71.428571428571
$ws.Range("M16").Value = -57.142857142857
$ws.Range("N16").Value = -95.488721804511
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = -50
$ws.Range("L17").Value = -55.555555555555
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -69.230769230769
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = 9.090909090909
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = -27.272727272727
$ws.Range("N18").Value = -90.839694656488
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -27.083333333333
$ws.Range("I19").Value = 102
$ws.Range("J19").Value = 147
$ws.Range("K19").Value = -30.612244897959
$ws.Range("L19").Value = -28.169014084507
$ws.Range("M19").Value = -32.894736842105
$ws.Range("N19").Value = -75.480769230769
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -100
$ws.Range("L20").Value = -54.545454545454
$ws.Range("M20").Value = 66.666666666666
$ws.Range("N20").Value = -96.575342465753
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -31.25
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -36.486486486486
$ws.Range("I21").Value = 150
$ws.Range("J21").Value = 219
$ws.Range("K21").Value = -31.506849315068
$ws.Range("L21").Value = -35.622317596566
$ws.Range("M21").Value = -30.232558139534
$ws.Range("N21").Value = -84.984984984985
$ws.Range("D22").Value = 2
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = -73.333333333333
$ws.Range("M22").Value = -63.636363636363
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 26.666666666666
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = -11.842105263157
$ws.Range("I24").Value = 220
$ws.Range("J24").Value = 219
$ws.Range("K24").Value = 0.456621004566
$ws.Range("L24").Value = 7.843137254901
$ws.Range("M24").Value = 48.648648648648
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = -34.375
$ws.Range("I25").Value = 162
$ws.Range("J25").Value = 181
$ws.Range("K25").Value = -10.49723756906
$ws.Range("L25").Value = 2.53164556962
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 49
$ws.Range("K26").Value = -38.775510204081
$ws.Range("L26").Value = -47.368421052631
$ws.Range("M26").Value = -44.444444444444
$ws.Range("G27").Value = 1
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -57.142857142857
$ws.Range("L28").Value = -40
